# Updated cryptos list on Sat Jun  1 20:42:57 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.767.85"
$ws.Range("D3").Value = "3.804.99"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.58"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.63"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.98"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "4.441.44"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "3.831.27"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "67.844.63"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.07"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.90"
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000152"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.42"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.10"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "3.953.44"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.46"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.53"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.07"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.749.55"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.43"
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.10"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.302"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.26"
$ws.Range("E46").Value = "  +11.90%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.95"
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.34"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.37"
$ws.Range("E49").Value = "  +9.25%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "147.73"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.85"
$ws.Range("E51").Value = "  +0.37%  "
